# Weekly update: add two new "Papa" price rows at the top of the data
# block (row 220) for the Terminal La Palmera de La Serena - Papa sheet,
# shifting all existing rows (220-269) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 220 (this pushes old rows 220..269
# down to 222..271 and grows the used range to A1:R271).
$ws.Rows.Item(220).Insert()
$ws.Rows.Item(220).Insert()

# --- New row 220 -----------------------------------------------------
$ws.Cells.Item(220, 1).Value = 8
$ws.Cells.Item(220, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(220, 3).Value = "Coquimbo"
$ws.Cells.Item(220, 4).Value = 44511
$ws.Cells.Item(220, 5).Value = 4
$ws.Cells.Item(220, 6).Value = 100114001
$ws.Cells.Item(220, 7).Value = "Papa"
$ws.Cells.Item(220, 8).Value = "Asterix"
$ws.Cells.Item(220, 9).Value = "1a nueva(o)"
$ws.Cells.Item(220, 10).Value = 2400
$ws.Cells.Item(220, 11).Value = 11500
$ws.Cells.Item(220, 12).Value = 12000
$ws.Cells.Item(220, 13).Value = 11750
$ws.Cells.Item(220, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(220, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(220, 16).Value = 470
$ws.Cells.Item(220, 17).Value = 25
$ws.Cells.Item(220, 18).Value = "Hortaliza"

# --- New row 221 -----------------------------------------------------
$ws.Cells.Item(221, 1).Value = 8
$ws.Cells.Item(221, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(221, 3).Value = "Coquimbo"
$ws.Cells.Item(221, 4).Value = 44511
$ws.Cells.Item(221, 5).Value = 4
$ws.Cells.Item(221, 6).Value = 100114001
$ws.Cells.Item(221, 7).Value = "Papa"
$ws.Cells.Item(221, 8).Value = "Cardinal"
$ws.Cells.Item(221, 9).Value = "1a (cosecha)"
$ws.Cells.Item(221, 10).Value = 2520
$ws.Cells.Item(221, 11).Value = 11500
$ws.Cells.Item(221, 12).Value = 12000
$ws.Cells.Item(221, 13).Value = 11750
$ws.Cells.Item(221, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(221, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(221, 16).Value = 470
$ws.Cells.Item(221, 17).Value = 25
$ws.Cells.Item(221, 18).Value = "Hortaliza"
